$d = $word.ActiveDocument

# Use Track Revisions so that Find/Replace creates a clean, isolated
# insert/delete revision pair for each edit; accepting that single
# revision (rather than a blanket AcceptAll) folds it back into plain
# runs without Word's "no-op retype" merging every other run in the
# paragraph (which is what produced the proofErr-wrapped run splits in
# the original document, and keeps the DMC condition line genuinely
# split into three runs as in the target).
$d.TrackRevisions = $true

# 1) Situación problemática paragraph: the only textual difference is
#    the removal of the spell-check run split around "promo" (and its
#    <w:proofErr> markers) — the visible text is unchanged, so simply
#    re-asserting the same text merges the three runs into one.
$old1 = "es que los necesita. Si vende exactamente la demanda mínima de cada promoción puede armar una cuarta promo. Debe alquilar un lugar para guardar los libros, y en caso de armar más de X promociones debe llamar a su sobrino para que lo ayude y le debe pagar."
$c1 = $d.Content
$c1.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $old1, 2) | Out-Null

# 2) "Precio de promo B)" paragraph: same kind of proofErr-run-split
#    cleanup, no actual text change.
$old2 = "Precio de promo B)"
$c2 = $d.Content
$c2.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $old2, 2) | Out-Null

# 3) DMC condition line: PB -> PC. Replace just the single "B"
#    character so the paragraph ends up split into exactly the three
#    runs shown in the target ("...<= P", "C", " <= DMC - 1 + D1").
$old3 = "DMC * D1 <= PB <= DMC - 1 + D1"
$c3 = $d.Content
$c3.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start = $c3.Start
$bChar = $d.Range($start + 13, $start + 14)
$bChar.Text = "C"

$d.TrackRevisions = $false

# Accept each revision individually instead of Revisions.AcceptAll():
# AcceptAll performs a whole-document normalization pass that strips
# unrelated w:rsid*/w:lastRenderedPageBreak bookkeeping from every run
# in the file, not just the ones we touched.
while ($d.Revisions.Count -gt 0) {
    $d.Revisions(1).Accept()
}
